$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values changed
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values changed / cleared
$ws.Range("B2").Value = 11.310938574386626
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = -0.17549645616645648
$ws.Range("E2").Value = -1.2083980962350438

# Row 3 values changed / cleared / added
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 5.6762230987349653
$ws.Range("D3").Value = 2.7539060664816475
$ws.Range("E3").Value = -4.2428102736428741

# Update the selection to match the edited range
$ws.Range("B1:E3").Select()
